# Revised "Preventative" to "Prevention" throughout the motivation column.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
[void]$used.Replace("Preventative", "Prevention", 1, 1, $true, $false, $false)

# Reset view/selection back to the top-left of the sheet (A1), matching
# the saved state in the target workbook (no scrolled/selected range).
[void]$ws.Range("A1").Select()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
